# Generate Report for Handoff
# Updates the localization status report: both zh-cn and de-de locales move
# from "In Translation" to "Ready for handoff", and the corresponding
# timestamps are refreshed to reflect the new handoff generation time.
# Also widens the date/status columns slightly to fit the new content.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-04 22:43:40"

# Widen the zh-cn / de-de status columns to fit "Ready for handoff"
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 22:43:35"
$wsZhCn.Range("C1").ColumnWidth = 16.33

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-04 22:43:40"
$wsDeDe.Range("C1").ColumnWidth = 16.33
